# Refresh the cryptocurrency Price (column D) and Volume(1h) (column E)
# figures pulled in by the scheduled GitHub Actions scraper job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''25.994.35'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").Value = '''1.641.52'
$ws.Range("E3").Value = '  +0.93%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = '''216.23'
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("E6").Value = '  +1.34%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +0.77%  '
$ws.Range("E9").Value = '  +1.30%  '
$ws.Range("E10").Value = '  +0.41%  '
$ws.Range("D11").Value = '''0.0795'
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("D12").Value = '''1.870.35'
$ws.Range("E12").Value = '  +0.68%  '
$ws.Range("E13").Value = '  +1.36%  '
$ws.Range("D14").Value = '''1.641.04'
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("E15").Value = '  +1.08%  '
$ws.Range("E16").Value = '  +1.43%  '
$ws.Range("D17").Value = '''62.96'
$ws.Range("E17").Value = '  +0.74%  '
$ws.Range("D18").Value = '''25.959.00'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = '''1.00'
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("D20").Value = '''193.24'
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").Value = '''9.94'
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").Value = '''6.26'
$ws.Range("E23").Value = '  +0.82%  '
$ws.Range("E24").Value = '  +7.85%  '
$ws.Range("E25").Value = '  +1.49%  '
$ws.Range("D26").Value = '''144.62'
$ws.Range("E26").Value = '  +1.75%  '
$ws.Range("E27").Value = '  +0.32%  '
$ws.Range("E28").Value = '  +1.48%  '
$ws.Range("E29").Value = '  +0.63%  '
$ws.Range("E30").Value = '  +0.49%  '
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("E33").Value = '  +1.41%  '
$ws.Range("E34").Value = '  -2.41%  '
$ws.Range("E35").Value = '  +2.60%  '
$ws.Range("E36").Value = '  +0.44%  '
$ws.Range("D37").Value = '''1.135.38'
$ws.Range("E37").Value = '  +0.69%  '
$ws.Range("E38").Value = '  -0.47%  '
$ws.Range("D39").Value = '''2.46'
$ws.Range("E39").Value = '  -0.20%  '
$ws.Range("E40").Value = '  +0.55%  '
$ws.Range("E41").Value = '  +1.20%  '
$ws.Range("D42").Value = '''99.44'
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("D43").Value = '''0.796'
$ws.Range("E43").Value = '  +0.32%  '
$ws.Range("D44").Value = '''1.779.60'
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("E45").Value = '  +2.73%  '
$ws.Range("D46").Value = '''56.72'
$ws.Range("E46").Value = '  +0.99%  '
$ws.Range("E47").Value = '  +2.39%  '
$ws.Range("E48").Value = '  +0.69%  '
$ws.Range("D49").Value = '''7.72'
$ws.Range("E49").Value = '  +1.45%  '
$ws.Range("D50").Value = '''0.415'
$ws.Range("E50").Value = '  +0.03%  '
$ws.Range("E51").Value = '  +0.57%  '
